# Auto-applied data refresh for Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1010.56525
$ws.Range("J32").Value = 1190.9166
$ws.Range("L32").Value = 1190.9166
$ws.Range("N32").Value = -1842.9166

$ws.Range("H41").Value = 306
$ws.Range("J41").Value = 499.5
$ws.Range("L41").Value = 499.5
$ws.Range("N41").Value = -1379.5

$ws.Range("H62").Value = 5498.5
$ws.Range("I62").Value = 4652.5
$ws.Range("K62").Value = 4652.5
$ws.Range("M62").Value = -4028.5

$ws.Range("H65").Value = 5498.5
$ws.Range("I65").Value = 4652.5
$ws.Range("K65").Value = 23262.5
$ws.Range("M65").Value = -20142.5

$ws.Range("H80").Value = 604.6667
$ws.Range("I80").Value = 899.75
$ws.Range("J80").Value = 368.6
$ws.Range("K80").Value = 2699.25
$ws.Range("L80").Value = 1105.8
$ws.Range("M80").Value = -1701.25
$ws.Range("N80").Value = -3101.8

$ws.Range("H83").Value = 604.6667
$ws.Range("I83").Value = 899.75
$ws.Range("J83").Value = 368.6
$ws.Range("K83").Value = 8097.75
$ws.Range("L83").Value = 3317.4
$ws.Range("M83").Value = -3105.75
$ws.Range("N83").Value = -13301.4

$ws.Range("H106").Value = 4698.1665
$ws.Range("I106").Value = 4698.1665
$ws.Range("K106").Value = 4698.1665
$ws.Range("M106").Value = -4067.1665

$ws.Range("H125").Value = 3959.8462
$ws.Range("J125").Value = 10632.667
$ws.Range("L125").Value = 95694.003
$ws.Range("N125").Value = -100614.003

$ws.Range("H137").Value = 3400.6
$ws.Range("I137").Value = 1999.6666
$ws.Range("K137").Value = 5998.9998
$ws.Range("M137").Value = -3448.9998

$ws.Range("H138").Value = 3432.9412
$ws.Range("I138").Value = 2070.375
$ws.Range("J138").Value = 4644.1113
$ws.Range("K138").Value = 6211.125
$ws.Range("L138").Value = 13932.3339
$ws.Range("M138").Value = -1071.125
$ws.Range("N138").Value = -24212.3339

$ws.Range("H141").Value = 2734.4211
$ws.Range("I141").Value = 2974.0588
$ws.Range("J141").Value = 697.5
$ws.Range("K141").Value = 8922.1764
$ws.Range("L141").Value = 2092.5
$ws.Range("M141").Value = -3742.1764
$ws.Range("N141").Value = -12452.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6121
$ws.Range("I32").Value = 6121
$ws.Range("K32").Value = 6121
$ws.Range("M32").Value = -5834

$ws.Range("H74").Value = 4877.4287
$ws.Range("I74").Value = 3976.111
$ws.Range("J74").Value = 6499.8
$ws.Range("K74").Value = 3976.111
$ws.Range("L74").Value = 6499.8
$ws.Range("M74").Value = -3102.111
$ws.Range("N74").Value = -8247.799999999999

$ws.Range("H77").Value = 4877.4287
$ws.Range("I77").Value = 3976.111
$ws.Range("J77").Value = 6499.8
$ws.Range("K77").Value = 19880.555
$ws.Range("L77").Value = 32499
$ws.Range("M77").Value = -15512.555
$ws.Range("N77").Value = -41235

$ws.Range("H97").Value = 1000.3333
$ws.Range("I97").Value = 961.3
$ws.Range("K97").Value = 961.3
$ws.Range("M97").Value = -465.3

$ws.Range("H132").Value = 1778.4117
$ws.Range("I132").Value = 1778.4117
$ws.Range("K132").Value = 5335.2351
$ws.Range("M132").Value = -2805.2351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3583.182
$ws.Range("I20").Value = 3003.3572
$ws.Range("J20").Value = 4597.875
$ws.Range("K20").Value = 3003.3572
$ws.Range("L20").Value = 4597.875
$ws.Range("M20").Value = -2756.3572
$ws.Range("N20").Value = -5091.875

$ws.Range("H80").Value = 899.5
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 899.5
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws.Range("H86").Value = 2499.3
$ws.Range("I86").Value = 2634.7144
$ws.Range("K86").Value = 2634.7144
$ws.Range("M86").Value = -1511.7144

$ws.Range("H89").Value = 2499.3
$ws.Range("I89").Value = 2634.7144
$ws.Range("K89").Value = 13173.572
$ws.Range("M89").Value = -7557.572

$ws.Range("H94").Value = 1746.5555
$ws.Range("I94").Value = 1818.32
$ws.Range("K94").Value = 1818.32
$ws.Range("M94").Value = -1367.32

$ws.Range("H105").Value = 1866.2307
$ws.Range("I105").Value = 1548.1
$ws.Range("K105").Value = 1548.1
$ws.Range("M105").Value = 198.9000000000001

$ws.Range("H107").Value = 1110.5
$ws.Range("I107").Value = 1110.5
$ws.Range("K107").Value = 1110.5
$ws.Range("M107").Value = 809.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2284.2183
$ws.Range("J31").Value = 3125.3333
$ws.Range("L31").Value = 3125.3333
$ws.Range("N31").Value = -3715.3333

$ws.Range("H34").Value = 2284.2183
$ws.Range("J34").Value = 3125.3333
$ws.Range("L34").Value = 3125.3333
$ws.Range("N34").Value = -3529.3333

$ws.Range("H41").Value = 16676.5
$ws.Range("J41").Value = 21250
$ws.Range("L41").Value = 21250
$ws.Range("N41").Value = -22106

$ws.Range("H50").Value = 28216.6

$ws.Range("H132").Value = 1940.2
$ws.Range("I132").Value = 1550.25
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 4650.75
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -2120.75
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3999.9
$ws.Range("J55").Value = 4333.222
$ws.Range("L55").Value = 12999.666
$ws.Range("N55").Value = -13353.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2659
$ws.Range("I80").Value = 1493.5
$ws.Range("J80").Value = 4990
$ws.Range("K80").Value = 1493.5
$ws.Range("L80").Value = 4990
$ws.Range("M80").Value = -495.5
$ws.Range("N80").Value = -6986

$ws.Range("H83").Value = 2659
$ws.Range("I83").Value = 1493.5
$ws.Range("J83").Value = 4990
$ws.Range("K83").Value = 7467.5
$ws.Range("L83").Value = 24950
$ws.Range("M83").Value = -2475.5
$ws.Range("N83").Value = -34934

$ws.Range("H97").Value = 661.6
$ws.Range("I97").Value = 446
$ws.Range("J97").Value = 1164.6666
$ws.Range("K97").Value = 446
$ws.Range("L97").Value = 1164.6666
$ws.Range("M97").Value = 50
$ws.Range("N97").Value = -2156.6666

$ws.Range("H102").Value = 4993
$ws.Range("I102").Value = 4993
$ws.Range("K102").Value = 4993
$ws.Range("M102").Value = -3371

$ws.Range("H122").Value = 3434
$ws.Range("I122").Value = 2223.8333
$ws.Range("K122").Value = 6671.499899999999
$ws.Range("M122").Value = -4221.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3061.9656
$ws.Range("I46").Value = 2253.4
$ws.Range("J46").Value = 3928.2856
$ws.Range("K46").Value = 2253.4
$ws.Range("L46").Value = 3928.2856
$ws.Range("M46").Value = -2065.4
$ws.Range("N46").Value = -4304.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 9102
$ws.Range("I113").Value = 17341.834
$ws.Range("J113").Value = 862.1667
$ws.Range("K113").Value = 52025.50199999999
$ws.Range("L113").Value = 2586.5001
$ws.Range("M113").Value = -49855.50199999999
$ws.Range("N113").Value = -6926.5001

$ws.Range("H132").Value = 1771.3572
$ws.Range("I132").Value = 1771.3572
$ws.Range("K132").Value = 5314.071599999999
$ws.Range("M132").Value = -2784.071599999999
